$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 127, shifting existing rows 127-199 down to 128-200.
$ws.Rows.Item(127).Insert()

# Populate the newly inserted row 127 with its data.
$ws.Cells.Item(127, 1).Value = 11
$ws.Cells.Item(127, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(127, 3).Value = "Bíobío"
$ws.Cells.Item(127, 4).Value = 44806
$ws.Cells.Item(127, 5).Value = 8
$ws.Cells.Item(127, 6).Value = "Fruta"
$ws.Cells.Item(127, 7).Value = 100108
$ws.Cells.Item(127, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(127, 9).Value = 100108005
$ws.Cells.Item(127, 10).Value = "Piña"
$ws.Cells.Item(127, 11).Value = "Caramelo"
$ws.Cells.Item(127, 12).Value = "Segunda"
$ws.Cells.Item(127, 13).Value = 200
$ws.Cells.Item(127, 14).Value = 18000
$ws.Cells.Item(127, 15).Value = 19000
$ws.Cells.Item(127, 16).Value = 18500
$ws.Cells.Item(127, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(127, 18).Value = "Ecuador"
$ws.Cells.Item(127, 19).Value = 1321
$ws.Cells.Item(127, 20).Value = 14

# Ensure the date cell keeps the same number format/style as the surrounding date column (style index 2).
$ws.Cells.Item(126, 4).Copy()
$ws.Cells.Item(127, 4).PasteSpecial(-4122)
